$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price/Volume columns for the data rows so that
# numeric-looking strings (e.g. "1.010", "0.06770") are preserved exactly as text
# instead of being auto-coerced into floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.426.45"
$ws.Range("E2").Value = "  -2.46%  "

$ws.Range("D3").Value = "1.937.14"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.95%  "

$ws.Range("D5").Value = "247.38"
$ws.Range("E5").Value = "  -2.44%  "

$ws.Range("D6").Value = "0.6899"
$ws.Range("E6").Value = "  -12.56%  "

$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +1.06%  "

$ws.Range("D8").Value = "0.3239"
$ws.Range("E8").Value = "  -5.02%  "

$ws.Range("D9").Value = "26.49"
$ws.Range("E9").Value = "  +3.44%  "

$ws.Range("D10").Value = "0.06770"
$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("D11").Value = "0.7933"
$ws.Range("E11").Value = "  -6.83%  "

$ws.Range("D12").Value = "0.08002"

$ws.Range("D13").Value = "1.949.90"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").Value = "5.378"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("D15").Value = "93.78"
$ws.Range("E15").Value = "  -8.26%  "

$ws.Range("D16").Value = "261.56"
$ws.Range("E16").Value = "  -6.21%  "

$ws.Range("D17").Value = "14.43"
$ws.Range("E17").Value = "  +3.41%  "

$ws.Range("D18").Value = "30.415.78"
$ws.Range("E18").Value = "  -2.37%  "

$ws.Range("D19").Value = "5.861"
$ws.Range("E19").Value = "  +3.16%  "

$ws.Range("D20").Value = "0.000007784"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("D21").Value = "2.214.12"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "1.009"
$ws.Range("E22").Value = "  +1.16%  "

$ws.Range("D23").Value = "1.011"
$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D24").Value = "6.821"
$ws.Range("E24").Value = "  +0.33%  "

$ws.Range("D25").Value = "9.608"
$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("D26").Value = "158.49"
$ws.Range("E26").Value = "  -4.09%  "

$ws.Range("D27").Value = "18.78"
$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("D28").Value = "2.249"
$ws.Range("E28").Value = "  +2.29%  "

$ws.Range("E29").Value = "  -16.21%  "

$ws.Range("D30").Value = "1.367"
$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("D31").Value = "1.556"
$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("D32").Value = "4.413"
$ws.Range("E32").Value = "  -3.05%  "

$ws.Range("D33").Value = "4.225"
$ws.Range("E33").Value = "  -2.73%  "

$ws.Range("D34").Value = "0.05075"
$ws.Range("E34").Value = "  -1.34%  "

$ws.Range("D35").Value = "1.189"
$ws.Range("E35").Value = "  -2.73%  "

$ws.Range("D36").Value = "0.7466"
$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("D37").Value = "2.733"
$ws.Range("E37").Value = "  -2.64%  "

$ws.Range("D38").Value = "0.01927"
$ws.Range("E38").Value = "  -3.06%  "

$ws.Range("D39").Value = "2.784"
$ws.Range("E39").Value = "  -4.32%  "

$ws.Range("D40").Value = "79.92"
$ws.Range("E40").Value = "  +1.84%  "

$ws.Range("D41").Value = "6.552"
$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("D42").Value = "2.042"
$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("D43").Value = "0.4419"
$ws.Range("E43").Value = "  -5.66%  "

$ws.Range("D44").Value = "1.009"
$ws.Range("E44").Value = "  +0.99%  "

$ws.Range("D45").Value = "0.8416"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("D46").Value = "101.41"
$ws.Range("E46").Value = "  -3.74%  "

$ws.Range("D47").Value = "9.771"
$ws.Range("E47").Value = "  -2.42%  "

$ws.Range("D48").Value = "7.287"
$ws.Range("E48").Value = "  -3.19%  "

$ws.Range("D49").Value = "35.93"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("D50").Value = "1.490"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05953"
$ws.Range("E51").Value = "  +0.67%  "
